$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text (row 1): username/password/res -> Username/Password/Result
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Result"

# Restyle the header row: drop the old green fill + border + big font,
# switch to a plain bold 11pt header, centered both ways with wrap text.
$headerRange = $ws.Range("A1:C1")
$headerRange.Interior.ColorIndex = -4142
$headerRange.Interior.Pattern = -4142
$headerRange.Borders.LineStyle = -4142
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

# Row 1 no longer needs the custom 21pt row height used by the old header style.
$ws.Rows.Item(1).AutoFit()

# Move the active selection (cosmetic, matches the saved cursor position).
$null = $ws.Range("F10").Select()
